# Cambiata gestione della colonna Anno in superdettagli
# Rename the worksheet "Sheet1" to "DoesNotMatter".
# The defined name "_xlnm._FilterDatabase" that points at Sheet1!$B$2:$J$46
# will automatically follow the rename since Excel keeps sheet-qualified
# references in defined names in sync with the sheet's name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "DoesNotMatter"
